# Updates cryptos list values (price/volume columns) to reflect the
# latest scrape, matching the commit "Updated cryptos list on Fri May 10
# 14:35:32 UTC 2024 with GitHub Actions".
#
# Rows 17/18 (Polkadot <-> WrappedBTC) and rows 39/40 (Stacks <-> OKB)
# swapped rank order, so their Coin/Link/Price/Volume values are fully
# rewritten rather than just the Price/Volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "62.102.39"
$ws.Cells.Item(2, 5).Value = "  +1.31%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.990.79"
$ws.Cells.Item(3, 5).Value = "  +0.34%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.19%  "

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "590.95"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.73%  "

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "150.21"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +4.44%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.35%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "2.980.39"
$ws.Cells.Item(8, 5).Value = "  -0.05%  "

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.506"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.26%  "

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.74"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +11.66%  "

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.147"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.67%  "

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.457"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +1.11%  "

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000228"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.85%  "

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "34.94"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.66%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  -0.48%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "3.493.61"
$ws.Cells.Item(16, 5).Value = "  +0.70%  "

# Row 17
$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "62.230.91"
$ws.Cells.Item(17, 5).Value = "  +1.60%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "Polkadot"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.98"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.12%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "2.998.14"
$ws.Cells.Item(19, 5).Value = "  +0.64%  "

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "441.72"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.74%  "

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "14.04"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.08%  "

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.686"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.81%  "

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.42"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.48%  "

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "82.05"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.42%  "

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.04"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +3.10%  "

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +2.65%  "

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.13"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +1.03%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.05%  "

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.39"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.90%  "

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.24"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +9.58%  "

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.68"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.47%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.26%  "

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "27.23"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.39%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.44%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "0.0₃0847"
$ws.Cells.Item(35, 5).Value = "  +5.33%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +1.50%  "

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.79"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.79%  "

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.05"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +8.17%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "OKB"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "50.11"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.07%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.06"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.59%  "

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.126"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.17%  "

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.90"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.57%  "

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "44.09"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +11.96%  "

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.299"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +11.83%  "

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0354"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.85%  "

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "379.40"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.84%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "2.682.41"
$ws.Cells.Item(47, 5).Value = "  +0.10%  "

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "132.86"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.55%  "

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "25.96"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +11.52%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +0.02%  "

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +3.97%  "
